# Update the "Status" sheet with refreshed daily figures and append new daily rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C (num of patients / "מספר מאומתים") for rows 42-96 ---
$colCUpdates = @{
    42 = 27
    43 = 31
    44 = 40
    45 = 44
    46 = 77
    47 = 96
    48 = 117
    49 = 152
    50 = 183
    51 = 221
    52 = 294
    53 = 389
    54 = 508
    55 = 645
    56 = 815
    57 = 1023
    58 = 1272
    59 = 1620
    60 = 2043
    61 = 2467
    62 = 2993
    63 = 3429
    64 = 3912
    65 = 4448
    66 = 4995
    67 = 5726
    68 = 6422
    69 = 7144
    70 = 7741
    71 = 8169
    72 = 8748
    73 = 9199
    74 = 9578
    75 = 9917
    76 = 10261
    77 = 10621
    78 = 10966
    79 = 11525
    80 = 11968
    81 = 12370
    82 = 12680
    83 = 12981
    84 = 13279
    85 = 13589
    86 = 13870
    87 = 14171
    88 = 14460
    89 = 14686
    90 = 14966
    91 = 15221
    92 = 15381
    93 = 15469
    94 = 15584
    95 = 15752
    96 = 15835
}
foreach ($row in $colCUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $colCUpdates[$row]
}

# --- Append new rows 97-101 with the latest data ---
$newRows = @(
    @{ Row = 97; A = 43951; B = 380339; C = 15979; D = 371; E = 110; F = 87; G = 223 }
    @{ Row = 98; A = 43952; B = 391104; C = 16096; D = 348; E = 108; F = 86; G = 227 }
    @{ Row = 99; A = 43953; B = 396659; C = 16153; D = 324; E = 109; F = 87; G = 230 }
    @{ Row = 100; A = 43954; B = 404586; C = 16182; D = 312; E = 94; F = 76; G = 233 }
    @{ Row = 101; A = 43955; B = 409515; C = 16237; D = 291; E = 91; F = 71; G = 236 }
)

# Existing date-style number format used by column A (matches cell style used for A2:A96)
$dateFormat = $ws.Range("A96").NumberFormat()

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}

Write-Output ("Final used range: " + $ws.UsedRange.Address())
